# Custom_Quota_Vary.xlsx - fix row index in the individual quota value
# placeholder ids (row-3 -> row-4) and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Columns I and K (rows 2-10) contained stale "row-3" placeholder ids;
# they should reference row 4 instead.
$ws.Range("I2:I10").Value = "individualQuotaValuesId-row-4-totalQuota"
$ws.Range("K2:K10").Value = "individualQuotaValuesId-row-4-colorQuota"

# Update the sheet's active selection/scroll position.
$ws.Range("K3:K10").Select()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
